$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.138.63'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.205.18'
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.83'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.53'
$ws.Range('E6').Value = '  -2.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.597'
$ws.Range('E7').Value = '  -4.58%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  -2.50%  '
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.392'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.764.55'
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '65.138.65'
$ws.Range('E14').Value = '  +0.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.65'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.202.71'
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000158'
$ws.Range('E17').Value = '  -0.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '412.59'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.92'
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.35'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.20'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.48'
$ws.Range('E23').Value = '  -1.32%  '
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.493'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('E26').Value = '  -5.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.92'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.85'
$ws.Range('E29').Value = '  -0.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.59'
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.97'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.42'
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '156.58'
$ws.Range('E34').Value = '  -0.87%  '
$ws.Range('E35').Value = '  -1.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.755.80'
$ws.Range('E36').Value = '  -2.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.73'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.33'
$ws.Range('E38').Value = '  -4.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.15'
$ws.Range('E39').Value = '  -1.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.715'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('E41').Value = '  +0.67%  '
$ws.Range('E42').Value = '  -1.89%  '
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '296.98'
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.57'
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0993'
$ws.Range('E46').Value = '  -1.37%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.96'
$ws.Range('E48').Value = '  -9.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.80'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.47'
$ws.Range('E50').Value = '  +0.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.903'
$ws.Range('E51').Value = '  -2.66%  '
